$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6.5
$ws.Range("F2").Value = 3

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 2

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 5
$ws.Range("F4").Value = 2

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 3
